$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.571.02"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "3.149.92"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'611.42"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("D6").Value = "'144.53"
$ws.Range("E6").Value = "  -2.25%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "3.147.65"
$ws.Range("E8").Value = "  -0.28%  "
$ws.Range("D9").Value = "'0.535"
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("D10").Value = "'0.152"
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("D11").Value = "'5.45"
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("D12").Value = "'0.478"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").Value = "'0.0000257"
$ws.Range("E13").Value = "  +2.04%  "
$ws.Range("D14").Value = "'35.71"
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("D15").Value = "3.664.33"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("E16").Value = "  +2.91%  "
$ws.Range("D17").Value = "64.505.64"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").Value = "3.144.72"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("D19").Value = "'6.91"
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("D20").Value = "'479.33"
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("D21").Value = "'14.81"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").Value = "'0.730"
$ws.Range("E22").Value = "  +2.51%  "
$ws.Range("D23").Value = "'7.88"
$ws.Range("E23").Value = "  +1.67%  "
$ws.Range("D24").Value = "'13.72"
$ws.Range("E24").Value = "  -0.36%  "
$ws.Range("D25").Value = "'85.58"
$ws.Range("E25").Value = "  +2.27%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "'8.60"
$ws.Range("E27").Value = "  +1.21%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "'2.79"
$ws.Range("E28").Value = "  -2.97%  "
$ws.Range("E29").Value = "  +8.55%  "
$ws.Range("D30").Value = "'2.10"
$ws.Range("E30").Value = "  -4.39%  "
$ws.Range("E31").Value = "  +3.25%  "
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("D33").Value = "'26.86"
$ws.Range("E33").Value = "  +2.19%  "
$ws.Range("D34").Value = "'2.66"
$ws.Range("E34").Value = "  -3.96%  "
$ws.Range("E35").Value = "  +1.20%  "
$ws.Range("D36").Value = "'5.99"
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0754"
$ws.Range("E37").Value = "  +4.02%  "
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").Value = "'52.85"
$ws.Range("E38").Value = "  -2.69%  "
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").Value = "'3.05"
$ws.Range("E39").Value = "  +5.28%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "'457.35"
$ws.Range("E40").Value = "  +1.09%  "
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("E42").Value = "  +1.28%  "
$ws.Range("D43").Value = "'8.37"
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("D44").Value = "2.882.44"
$ws.Range("E44").Value = "  +1.00%  "
$ws.Range("D45").Value = "'0.265"
$ws.Range("E45").Value = "  -0.97%  "
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("E47").Value = "  +5.86%  "
$ws.Range("D48").Value = "'26.67"
$ws.Range("E48").Value = "  +0.85%  "
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("D50").Value = "'0.115"
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("D51").Value = "'121.35"
$ws.Range("E51").Value = "  +1.61%  "
